$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 50 gets bumped from date serial
# 46081 (2026-02-28) to 46082 (2026-03-01).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 3).Value = 46082
}
